# Update cryptocurrency price and 1h volume/change values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '26.434.68'
Set-TextValue $ws.Range('D3') '1.811.20'
Set-TextValue $ws.Range('E3') '  +0.32%  '
Set-TextValue $ws.Range('D4') '1.004'
Set-TextValue $ws.Range('E4') '  -0.57%  '
Set-TextValue $ws.Range('E5') '  -0.51%  '
Set-TextValue $ws.Range('D6') '305.37'
Set-TextValue $ws.Range('E6') '  -0.94%  '
Set-TextValue $ws.Range('E7') '  -0.74%  '
Set-TextValue $ws.Range('D8') '0.3583'
Set-TextValue $ws.Range('D9') '46.31'
Set-TextValue $ws.Range('E9') '  +2.74%  '
Set-TextValue $ws.Range('D10') '0.07056'
Set-TextValue $ws.Range('E10') '  -0.84%  '
Set-TextValue $ws.Range('D11') '0.8887'
Set-TextValue $ws.Range('E11') '  +1.55%  '
Set-TextValue $ws.Range('D12') '0.07780'
Set-TextValue $ws.Range('E12') '  +0.51%  '
Set-TextValue $ws.Range('D13') '19.29'
Set-TextValue $ws.Range('E13') '  -0.20%  '
Set-TextValue $ws.Range('D14') '1.791.45'
Set-TextValue $ws.Range('E14') '  -1.38%  '
Set-TextValue $ws.Range('E15') '  -0.01%  '
Set-TextValue $ws.Range('D16') '6.303'
Set-TextValue $ws.Range('E16') '  -0.58%  '
Set-TextValue $ws.Range('D17') '84.87'
Set-TextValue $ws.Range('E17') '  -1.24%  '
Set-TextValue $ws.Range('E18') '  -0.53%  '
Set-TextValue $ws.Range('D19') '0.000008517'
Set-TextValue $ws.Range('E19') '  -0.61%  '
Set-TextValue $ws.Range('E20') '  -0.56%  '
Set-TextValue $ws.Range('D21') '26.483.07'
Set-TextValue $ws.Range('E21') '  -0.17%  '
Set-TextValue $ws.Range('D22') '14.13'
Set-TextValue $ws.Range('E22') '  -0.70%  '
Set-TextValue $ws.Range('D23') '4.953'
Set-TextValue $ws.Range('E23') '  -0.34%  '
Set-TextValue $ws.Range('D24') '2.023.40'
Set-TextValue $ws.Range('E24') '  -1.02%  '
Set-TextValue $ws.Range('E25') '  +0.83%  '
Set-TextValue $ws.Range('E26') '  -0.93%  '
Set-TextValue $ws.Range('D27') '150.91'
Set-TextValue $ws.Range('E27') '  +0.20%  '
Set-TextValue $ws.Range('D28') '17.77'
Set-TextValue $ws.Range('E28') '  -0.78%  '
Set-TextValue $ws.Range('D29') '2.046'
Set-TextValue $ws.Range('E29') '  +2.56%  '
Set-TextValue $ws.Range('D30') '112.22'
Set-TextValue $ws.Range('E30') '  -0.31%  '
Set-TextValue $ws.Range('D31') '4.822'
Set-TextValue $ws.Range('E31') '  +0.09%  '
Set-TextValue $ws.Range('D32') '0.08679'
Set-TextValue $ws.Range('D33') '3.139'
Set-TextValue $ws.Range('E33') '  +2.99%  '
Set-TextValue $ws.Range('D34') '0.7446'
Set-TextValue $ws.Range('E34') '  +2.31%  '
Set-TextValue $ws.Range('D35') '2.725'
Set-TextValue $ws.Range('E35') '  +6.33%  '
Set-TextValue $ws.Range('D36') '4.426'
Set-TextValue $ws.Range('E36') '  +0.00%  '
Set-TextValue $ws.Range('E37') '  -0.38%  '
Set-TextValue $ws.Range('E38') '  -1.23%  '
Set-TextValue $ws.Range('D39') '0.01925'
Set-TextValue $ws.Range('E39') '  -0.16%  '
Set-TextValue $ws.Range('D40') '2.902'
Set-TextValue $ws.Range('E40') '  +0.75%  '
Set-TextValue $ws.Range('D41') '0.05088'
Set-TextValue $ws.Range('D42') '0.5073'
Set-TextValue $ws.Range('E42') '  +1.65%  '
Set-TextValue $ws.Range('D43') '6.742'
Set-TextValue $ws.Range('E43') '  -2.77%  '
Set-TextValue $ws.Range('D44') '0.1505'
Set-TextValue $ws.Range('E44') '  -3.70%  '
Set-TextValue $ws.Range('D45') '8.036'
Set-TextValue $ws.Range('E45') '  -0.80%  '
Set-TextValue $ws.Range('D46') '0.4701'
Set-TextValue $ws.Range('E46') '  +2.42%  '
Set-TextValue $ws.Range('D47') '1.003'
Set-TextValue $ws.Range('E47') '  -0.61%  '
Set-TextValue $ws.Range('D48') '9.997'
Set-TextValue $ws.Range('E48') '  +0.71%  '
Set-TextValue $ws.Range('D49') '100.08'
Set-TextValue $ws.Range('E49') '  -1.59%  '
Set-TextValue $ws.Range('E50') '  -0.79%  '
Set-TextValue $ws.Range('D51') '0.05986'
Set-TextValue $ws.Range('E51') '  -0.06%  '
